# Apply the two citation/page-count fixes described by the diff:
#  1. Footer "Page X of Y": the cached result of the PAGE field (X) should
#     read "5" instead of "1" (document grew to 5 pages).
#  2. Header citation: "City & Community 18(1), 280." should become
#     "City & Community, 18(1), 280-301." -- i.e. add a comma right after
#     the (italic) journal title and extend the page number to a full
#     page range by appending "-301".

$d = $word.ActiveDocument
$sec = $d.Sections.First

# --- 1. Footer: fix the PAGE field's cached/display text -------------------
$ftr = $sec.Footers.Item(1)
$pageField = $ftr.Range.Fields.Item(1)          # the " PAGE \* Arabic ..." field
$pageResult = $pageField.Result                  # Range covering just its result text ("1")
$null = $pageResult.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "5", 2)

# --- 2. Header: fix the journal citation punctuation ------------------------
$hdr = $sec.Headers.Item(1)
$cite = $hdr.Range.Duplicate
$null = $cite.Find.Execute(" 18(1), 280", $true, $false, $false, $false, $false, $true, 1, $false, ", 18(1), 280-301", 2)
